$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.401.07'
$ws.Range("E2").Value = '  +4.33%  '
$ws.Range("D3").Value = '3.487.04'
$ws.Range("E3").Value = '  +3.60%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.53%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.478'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.48%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.72'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.20%  '
$ws.Range("E11").Value = '  +4.25%  '
$ws.Range("D12").Value = '4.084.48'
$ws.Range("E12").Value = '  +3.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +6.88%  '
$ws.Range("E14").Value = '  -0.23%  '
$ws.Range("D15").Value = '3.491.91'
$ws.Range("E15").Value = '  +4.04%  '
$ws.Range("E16").Value = '  +4.03%  '
$ws.Range("D17").Value = '63.385.28'
$ws.Range("E17").Value = '  +4.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.13%  '
$ws.Range("E19").Value = '  +6.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.40'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '393.12'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.90%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.565'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.21'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  +8.49%  '
$ws.Range("D26").Value = '3.628.43'
$ws.Range("E26").Value = '  +3.66%  '
$ws.Range("E27").Value = '  -2.70%  '
$ws.Range("E28").Value = '  +10.43%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("E30").Value = '  +5.47%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.79%  '
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.15'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.38%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("E34").Value = '  +3.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '32.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +27.22%  '
$ws.Range("E36").Value = '  +8.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.15'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.16%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '171.46'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.57'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.63%  '
$ws.Range("D40").Value = '3.525.25'
$ws.Range("E40").Value = '  +3.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0767'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.81%  '
$ws.Range("E42").Value = '  +4.68%  '
$ws.Range("E43").Value = '  +4.17%  '
$ws.Range("E44").Value = '  +7.55%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '42.52'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.70%  '
$ws.Range("E46").Value = '  +10.43%  '
$ws.Range("D47").Value = '2.622.24'
$ws.Range("E47").Value = '  +7.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.87'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.18%  '
$ws.Range("E49").Value = '  +18.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.75'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.17%  '
$ws.Range("E51").Value = '  +5.22%  '
